# Import template: add four new fields (color, solubility, form, inventory
# label) as additional columns, inserted right before the existing
# "molfile" / "canonical smiles" columns, and update the sheet's view
# state to reflect where the author was working when they made the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank columns at X:AA. This pushes the existing "molfile"
# (previously X) and "canonical smiles" (previously Y) columns two slots
# to the right (to AB/AC respectively), carrying their values/styles with
# them, and widens the used range from A1:Y6 to A1:AC6 automatically.
$null = $ws.Range("X:AA").Insert()

# Populate the headers of the newly inserted columns.
$ws.Range("X1").Value = "color"
$ws.Range("Y1").Value = "solubility"
$ws.Range("Z1").Value = "form"
$ws.Range("AA1").Value = "inventory label"

# "inventory label" needs a wider column, matching the author's resize.
$ws.Columns("AA").ColumnWidth = 17.5

# Reflect the cell selection / scroll position left by the author.
$null = $ws.Range("Z5").Select()
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
